$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume data columns keep their text formatting so Excel
# does not silently reinterpret numeric-looking strings as numbers.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '64.498.65'
$ws.Range('E2').Value = '  -1.64%  '
$ws.Range('D3').Value = '2.633.90'
$ws.Range('E3').Value = '  -1.14%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '580.26'
$ws.Range('E5').Value = '  -2.89%  '
$ws.Range('D6').Value = '156.73'
$ws.Range('E6').Value = '  +0.14%  '
$ws.Range('D7').Value = '0.652'
$ws.Range('E7').Value = '  +5.97%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -3.64%  '
$ws.Range('D10').Value = '5.82'
$ws.Range('E10').Value = '  +0.24%  '
$ws.Range('E11').Value = '  -1.93%  '
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = '28.69'
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('E14').Value = '  -5.42%  '
$ws.Range('D15').Value = '3.107.86'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').Value = '64.294.23'
$ws.Range('E16').Value = '  -1.80%  '
$ws.Range('D17').Value = '2.634.31'
$ws.Range('E17').Value = '  -0.98%  '
$ws.Range('D18').Value = '12.28'
$ws.Range('E18').Value = '  -2.83%  '
$ws.Range('D19').Value = '4.68'
$ws.Range('E19').Value = '  -2.03%  '
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('D21').Value = '346.85'
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = '68.06'
$ws.Range('E23').Value = '  -1.46%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('E25').Value = '  +4.50%  '
$ws.Range('D26').Value = '9.41'
$ws.Range('E26').Value = '  -3.07%  '

# Rows 27 and 28 swapped order (Bittensor now ranked above Fetch.AI)
$ws.Range('B27').Value = 'Bittensor'
$ws.Range('C27').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D27').Value = '582.34'
$ws.Range('E27').Value = '  +9.95%  '
$ws.Range('B28').Value = 'Fetch.AI'
$ws.Range('C28').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D28').Value = '1.57'
$ws.Range('E28').Value = '  -1.06%  '

$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('D31').Value = '7.93'
$ws.Range('E31').Value = '  -0.47%  '
$ws.Range('E32').Value = '  -2.25%  '
$ws.Range('D33').Value = '6.69'
$ws.Range('E33').Value = '  +3.98%  '
$ws.Range('E34').Value = '  -3.25%  '
$ws.Range('E35').Value = '  -2.25%  '
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('D37').Value = '20.11'
$ws.Range('E37').Value = '  -2.19%  '
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('D40').Value = '154.95'
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('E42').Value = '  +6.31%  '
$ws.Range('D43').Value = '158.24'
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('D44').Value = '4.00'
$ws.Range('E44').Value = '  -2.07%  '
$ws.Range('D45').Value = '23.27'
$ws.Range('E45').Value = '  +3.22%  '
$ws.Range('D46').Value = '0.0600'
$ws.Range('E46').Value = '  -1.24%  '
$ws.Range('D47').Value = '0.637'
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('E48').Value = '  +4.06%  '
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('E50').Value = '  -3.10%  '
$ws.Range('E51').Value = '  -5.71%  '
